# Update "想去人数" (want-to-go count) figures for three events that appear
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F5").Value = 455
$wsExpo.Range("F7").Value = 2480
$wsExpo.Range("F9").Value = 6546

# Sheet 4: 全部类型 (All types)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F5").Value = 455
$wsAll.Range("F9").Value = 2480
$wsAll.Range("F11").Value = 6546
